# Append the new Adafruit IO reading as row 31, matching the existing
# "Timestamp | Feed Key | Value | Latitude | Longitude | Elevation" layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 31

$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"

# Column C holds "25" as text in the source data (all cells in this sheet
# are stored as text, including numeric-looking ones), so force a text
# number format before assigning to stop auto-conversion to a number.
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "25"

$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"
